$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.415.38"
$ws.Range("E2").Value = "  -1.96%  "
$ws.Range("D3").Value = "2.982.24"
$ws.Range("E3").Value = "  -1.36%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.12"
$ws.Range("E5").Value = "  +1.91%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.44"
$ws.Range("E6").Value = "  -2.52%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").Value = "2.982.05"
$ws.Range("E8").Value = "  -1.29%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.512"
$ws.Range("E9").Value = "  -2.90%  "
$ws.Range("E10").Value = "  -1.65%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.02"
$ws.Range("E11").Value = "  +2.73%  "
$ws.Range("E12").Value = "  -0.93%  "
$ws.Range("E13").Value = "  -1.70%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.07"
$ws.Range("E14").Value = "  -2.15%  "
$ws.Range("E15").Value = "  +3.46%  "
$ws.Range("D16").Value = "3.477.63"
$ws.Range("E16").Value = "  -1.26%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.90"
$ws.Range("E17").Value = "  -3.37%  "
$ws.Range("D18").Value = "61.378.29"
$ws.Range("E18").Value = "  -1.96%  "
$ws.Range("D19").Value = "2.983.26"
$ws.Range("E19").Value = "  -1.26%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "444.72"
$ws.Range("E20").Value = "  -4.33%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.89"
$ws.Range("E22").Value = "  -1.24%  "
$ws.Range("E23").Value = "  -2.36%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "81.06"
$ws.Range("E24").Value = "  -0.87%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.78"
$ws.Range("E25").Value = "  +3.97%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.17"
$ws.Range("E26").Value = "  -4.61%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.00"
$ws.Range("E27").Value = "  -3.46%  "
$ws.Range("E28").Value = "  +0.11%  "
$ws.Range("E29").Value = "  +2.22%  "
$ws.Range("E30").Value = "  -0.08%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.24"
$ws.Range("E31").Value = "  +1.33%  "
$ws.Range("E32").Value = "  -3.69%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.10"
$ws.Range("E33").Value = "  -5.99%  "
$ws.Range("E34").Value = "  -0.08%  "
$ws.Range("D35").Value = "0.0₃0808"
$ws.Range("E35").Value = "  -0.12%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.01"
$ws.Range("E36").Value = "  -1.20%  "
$ws.Range("E37").Value = "  -0.92%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "50.20"
$ws.Range("E38").Value = "  -0.58%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.93"
$ws.Range("E39").Value = "  -2.77%  "
$ws.Range("E40").Value = "  -5.80%  "
$ws.Range("E41").Value = "  +7.36%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.85"
$ws.Range("E42").Value = "  -3.86%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "388.49"
$ws.Range("E43").Value = "  -2.79%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "39.47"
$ws.Range("E44").Value = "  +5.47%  "
$ws.Range("E45").Value = "  -3.72%  "
$ws.Range("E46").Value = "  -3.36%  "
$ws.Range("D47").Value = "2.682.47"
$ws.Range("E47").Value = "  -2.68%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "132.00"
$ws.Range("E48").Value = "  +2.31%  "
$ws.Range("E49").Value = "  +0.14%  "
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.107"
$ws.Range("E50").Value = "  -2.72%  "
$ws.Range("B51").Value = "ThetaToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.15"
$ws.Range("E51").Value = "  -2.65%  "
